$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 61273.445  # H11: 68933.125 -> 61273.445
$ws.Cells.Item(11, 9).Value = 61273.445  # I11: 68933.125 -> 61273.445
$ws.Cells.Item(11, 11).Value = 61273.445  # K11: 68933.125 -> 61273.445
$ws.Cells.Item(11, 13).Value = -61133.445  # M11: -68793.125 -> -61133.445

$ws.Cells.Item(42, 8).Value = 80.5  # H42: 115.71429 -> 80.5
$ws.Cells.Item(42, 9).Value = 8  # I42: 167.5 -> 8
$ws.Cells.Item(42, 11).Value = 24  # K42: 502.5 -> 24
$ws.Cells.Item(42, 13).Value = 206  # M42: -272.5 -> 206

$ws.Cells.Item(43, 8).Value = 1786.75  # H43: 1211.75 -> 1786.75
$ws.Cells.Item(43, 9).Value = 2349  # I43: 1199 -> 2349
$ws.Cells.Item(43, 11).Value = 2349  # K43: 1199 -> 2349
$ws.Cells.Item(43, 13).Value = -2280  # M43: -1130 -> -2280

$ws.Cells.Item(53, 8).Value = 2767.6956  # H53: 2769 -> 2767.6956
$ws.Cells.Item(53, 9).Value = 75.583336  # I53: 88.09999999999999 -> 75.583336
$ws.Cells.Item(53, 10).Value = 5704.5454  # J53: 4831.231 -> 5704.5454
$ws.Cells.Item(53, 11).Value = 75.583336  # K53: 88.09999999999999 -> 75.583336
$ws.Cells.Item(53, 12).Value = 5704.5454  # L53: 4831.231 -> 5704.5454
$ws.Cells.Item(53, 13).Value = 561.416664  # M53: 548.9 -> 561.416664
$ws.Cells.Item(53, 14).Value = -6978.5454  # N53: -6105.231 -> -6978.5454

$ws.Cells.Item(62, 8).Value = 1205278.2  # H62: 120000 -> 1205278.2
$ws.Cells.Item(62, 9).Value = 3435834.8  # I62: 300000 -> 3435834.8
$ws.Cells.Item(62, 11).Value = 3435834.8  # K62: 300000 -> 3435834.8
$ws.Cells.Item(62, 13).Value = -3435210.8  # M62: -299376 -> -3435210.8

$ws.Cells.Item(65, 8).Value = 1205278.2  # H65: 120000 -> 1205278.2
$ws.Cells.Item(65, 9).Value = 3435834.8  # I65: 300000 -> 3435834.8
$ws.Cells.Item(65, 11).Value = 17179174  # K65: 1500000 -> 17179174
$ws.Cells.Item(65, 13).Value = -17176054  # M65: -1496880 -> -17176054

$ws.Cells.Item(86, 8).Value = 8721465  # H86: 8721531 -> 8721465
$ws.Cells.Item(86, 9).Value = 3885.9092  # I86: 4125.5 -> 3885.9092
$ws.Cells.Item(86, 10).Value = 16712580  # J86: 15427227 -> 16712580
$ws.Cells.Item(86, 11).Value = 3885.9092  # K86: 4125.5 -> 3885.9092
$ws.Cells.Item(86, 12).Value = 16712580  # L86: 15427227 -> 16712580
$ws.Cells.Item(86, 13).Value = -2762.9092  # M86: -3002.5 -> -2762.9092
$ws.Cells.Item(86, 14).Value = -16714826  # N86: -15429473 -> -16714826

$ws.Cells.Item(88, 8).Value = 1238.5333  # H88: 1149.5217 -> 1238.5333
$ws.Cells.Item(88, 10).Value = 1348.25  # J88: 1165.4375 -> 1348.25
$ws.Cells.Item(88, 12).Value = 1348.25  # L88: 1165.4375 -> 1348.25
$ws.Cells.Item(88, 14).Value = -2160.25  # N88: -1977.4375 -> -2160.25

$ws.Cells.Item(89, 8).Value = 8721465  # H89: 8721531 -> 8721465
$ws.Cells.Item(89, 9).Value = 3885.9092  # I89: 4125.5 -> 3885.9092
$ws.Cells.Item(89, 10).Value = 16712580  # J89: 15427227 -> 16712580
$ws.Cells.Item(89, 11).Value = 19429.546  # K89: 20627.5 -> 19429.546
$ws.Cells.Item(89, 12).Value = 83562900  # L89: 77136135 -> 83562900
$ws.Cells.Item(89, 13).Value = -13813.546  # M89: -15011.5 -> -13813.546
$ws.Cells.Item(89, 14).Value = -83574132  # N89: -77147367 -> -83574132

$ws.Cells.Item(91, 8).Value = 1238.5333  # H91: 1149.5217 -> 1238.5333
$ws.Cells.Item(91, 10).Value = 1348.25  # J91: 1165.4375 -> 1348.25
$ws.Cells.Item(91, 12).Value = 1348.25  # L91: 1165.4375 -> 1348.25
$ws.Cells.Item(91, 14).Value = -4156.25  # N91: -3973.4375 -> -4156.25

$ws.Cells.Item(92, 8).Value = 45696.375  # H92: 47682.523 -> 45696.375
$ws.Cells.Item(92, 10).Value = 136216.12  # J92: 155673.42 -> 136216.12
$ws.Cells.Item(92, 12).Value = 136216.12  # L92: 155673.42 -> 136216.12
$ws.Cells.Item(92, 14).Value = -138712.12  # N92: -158169.42 -> -138712.12

$ws.Cells.Item(116, 8).Value = 44958010  # H116: 34379940 -> 44958010
$ws.Cells.Item(116, 9).Value = 41840500  # I116: 25104788 -> 41840500
$ws.Cells.Item(116, 11).Value = 41840500  # K116: 25104788 -> 41840500
$ws.Cells.Item(116, 13).Value = -41837058  # M116: -25101346 -> -41837058

$ws.Cells.Item(132, 8).Value = 2809.256  # H132: 2839.926 -> 2809.256
$ws.Cells.Item(132, 9).Value = 2628.1018  # I132: 2667.8103 -> 2628.1018
$ws.Cells.Item(132, 11).Value = 7884.305399999999  # K132: 8003.4309 -> 7884.305399999999
$ws.Cells.Item(132, 13).Value = -5354.305399999999  # M132: -5473.4309 -> -5354.305399999999

$ws.Cells.Item(137, 8).Value = 2158.432  # H137: 2158.955 -> 2158.432
$ws.Cells.Item(137, 9).Value = 1578.6818  # I137: 1521.8334 -> 1578.6818
$ws.Cells.Item(137, 10).Value = 2351.682  # J137: 2394.2 -> 2351.682
$ws.Cells.Item(137, 11).Value = 4736.0454  # K137: 4565.5002 -> 4736.0454
$ws.Cells.Item(137, 12).Value = 7055.045999999999  # L137: 7182.599999999999 -> 7055.045999999999
$ws.Cells.Item(137, 13).Value = -2186.0454  # M137: -2015.5002 -> -2186.0454
$ws.Cells.Item(137, 14).Value = -12155.046  # N137: -12282.6 -> -12155.046

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(21, 8).Value = 9254.25  # H21: 4288.857 -> 9254.25
$ws.Cells.Item(21, 9).Value = 7000  # I21: 3080.3076 -> 7000
$ws.Cells.Item(21, 10).Value = 11508.5  # J21: 20000 -> 11508.5
$ws.Cells.Item(21, 11).Value = 7000  # K21: 3080.3076 -> 7000
$ws.Cells.Item(21, 12).Value = 11508.5  # L21: 20000 -> 11508.5
$ws.Cells.Item(21, 13).Value = -6626  # M21: -2706.3076 -> -6626
$ws.Cells.Item(21, 14).Value = -12256.5  # N21: -20748 -> -12256.5

$ws.Cells.Item(32, 8).Value = 265569.84  # H32: 247935.48 -> 265569.84
$ws.Cells.Item(32, 9).Value = 280118.16  # I32: 325176.22 -> 280118.16
$ws.Cells.Item(32, 10).Value = 3700  # J32: 8489.200000000001 -> 3700
$ws.Cells.Item(32, 11).Value = 280118.16  # K32: 325176.22 -> 280118.16
$ws.Cells.Item(32, 12).Value = 3700  # L32: 8489.200000000001 -> 3700
$ws.Cells.Item(32, 13).Value = -279831.16  # M32: -324889.22 -> -279831.16
$ws.Cells.Item(32, 14).Value = -4274  # N32: -9063.200000000001 -> -4274

$ws.Cells.Item(45, 8).Value = 2983.7273  # H45: 2574.0667 -> 2983.7273
$ws.Cells.Item(45, 9).Value = 2521.4  # I45: 2044.1111 -> 2521.4
$ws.Cells.Item(45, 11).Value = 2521.4  # K45: 2044.1111 -> 2521.4
$ws.Cells.Item(45, 13).Value = -2144.4  # M45: -1667.1111 -> -2144.4

$ws.Cells.Item(60, 8).Value = 11371.333  # H60: 15374.5 -> 11371.333
$ws.Cells.Item(60, 9).Value = 18000  # I60: 15374.5 -> 18000
$ws.Cells.Item(60, 10).Value = 8057  # J60: 0 -> 8057
$ws.Cells.Item(60, 11).Value = 18000  # K60: 15374.5 -> 18000
$ws.Cells.Item(60, 12).Value = 8057  # L60: 0 -> 8057
$ws.Cells.Item(60, 13).Value = -17267  # M60: -14641.5 -> -17267
$ws.Cells.Item(60, 14).Value = -9523  # N60: None -> -9523

$ws.Cells.Item(122, 8).Value = 3116.2  # H122: 3203.9473 -> 3116.2
$ws.Cells.Item(122, 9).Value = 2027.4  # I122: 2091.6667 -> 2027.4
$ws.Cells.Item(122, 11).Value = 6082.200000000001  # K122: 6275.000100000001 -> 6082.200000000001
$ws.Cells.Item(122, 13).Value = -3632.200000000001  # M122: -3825.000100000001 -> -3632.200000000001

$ws.Cells.Item(132, 8).Value = 626875.5  # H132: 647060.25 -> 626875.5
$ws.Cells.Item(132, 9).Value = 386212.2  # I132: 401614.7 -> 386212.2
$ws.Cells.Item(132, 11).Value = 1158636.6  # K132: 1204844.1 -> 1158636.6
$ws.Cells.Item(132, 13).Value = -1156106.6  # M132: -1202314.1 -> -1156106.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2604.4443  # H134: 2766.75 -> 2604.4443
$ws.Cells.Item(134, 9).Value = 2123.3076  # I134: 2271.9092 -> 2123.3076
$ws.Cells.Item(134, 11).Value = 6369.9228  # K134: 6815.7276 -> 6369.9228
$ws.Cells.Item(134, 13).Value = -3834.9228  # M134: -4280.7276 -> -3834.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2467.602  # H31: 2094.71 -> 2467.602
$ws.Cells.Item(31, 9).Value = 1135.6428  # I31: 857.86487 -> 1135.6428
$ws.Cells.Item(31, 10).Value = 3000.3857  # J31: 2821.111 -> 3000.3857
$ws.Cells.Item(31, 11).Value = 1135.6428  # K31: 857.86487 -> 1135.6428
$ws.Cells.Item(31, 12).Value = 3000.3857  # L31: 2821.111 -> 3000.3857
$ws.Cells.Item(31, 13).Value = -840.6428000000001  # M31: -562.86487 -> -840.6428000000001
$ws.Cells.Item(31, 14).Value = -3590.3857  # N31: -3411.111 -> -3590.3857

$ws.Cells.Item(34, 8).Value = 2467.602  # H34: 2094.71 -> 2467.602
$ws.Cells.Item(34, 9).Value = 1135.6428  # I34: 857.86487 -> 1135.6428
$ws.Cells.Item(34, 10).Value = 3000.3857  # J34: 2821.111 -> 3000.3857
$ws.Cells.Item(34, 11).Value = 1135.6428  # K34: 857.86487 -> 1135.6428
$ws.Cells.Item(34, 12).Value = 3000.3857  # L34: 2821.111 -> 3000.3857
$ws.Cells.Item(34, 13).Value = -933.6428000000001  # M34: -655.86487 -> -933.6428000000001
$ws.Cells.Item(34, 14).Value = -3404.3857  # N34: -3225.111 -> -3404.3857

$ws.Cells.Item(39, 8).Value = 925.1667  # H39: 999.6667 -> 925.1667
$ws.Cells.Item(39, 9).Value = 925.1667  # I39: 999.6667 -> 925.1667
$ws.Cells.Item(39, 11).Value = 925.1667  # K39: 999.6667 -> 925.1667
$ws.Cells.Item(39, 13).Value = -534.1667  # M39: -608.6667 -> -534.1667

$ws.Cells.Item(49, 8).Value = 925.1667  # H49: 999.6667 -> 925.1667
$ws.Cells.Item(49, 9).Value = 925.1667  # I49: 999.6667 -> 925.1667
$ws.Cells.Item(49, 11).Value = 925.1667  # K49: 999.6667 -> 925.1667
$ws.Cells.Item(49, 13).Value = -743.1667  # M49: -817.6667 -> -743.1667

$ws.Cells.Item(134, 8).Value = 2346.1082  # H134: 2374.7896 -> 2346.1082
$ws.Cells.Item(134, 9).Value = 1969.9524  # I134: 1966.9048 -> 1969.9524
$ws.Cells.Item(134, 10).Value = 2839.8125  # J134: 2878.647 -> 2839.8125
$ws.Cells.Item(134, 11).Value = 5909.857199999999  # K134: 5900.7144 -> 5909.857199999999
$ws.Cells.Item(134, 12).Value = 8519.4375  # L134: 8635.940999999999 -> 8519.4375
$ws.Cells.Item(134, 13).Value = -3374.857199999999  # M134: -3365.7144 -> -3374.857199999999
$ws.Cells.Item(134, 14).Value = -13589.4375  # N134: -13705.941 -> -13589.4375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1806.9524  # H68: 1827.3 -> 1806.9524
$ws.Cells.Item(68, 10).Value = 2018.5  # J68: 2059.7334 -> 2018.5
$ws.Cells.Item(68, 12).Value = 6055.5  # L68: 6179.2002 -> 6055.5
$ws.Cells.Item(68, 14).Value = -7677.5  # N68: -7801.2002 -> -7677.5

$ws.Cells.Item(71, 8).Value = 1806.9524  # H71: 1827.3 -> 1806.9524
$ws.Cells.Item(71, 10).Value = 2018.5  # J71: 2059.7334 -> 2018.5
$ws.Cells.Item(71, 12).Value = 18166.5  # L71: 18537.6006 -> 18166.5
$ws.Cells.Item(71, 14).Value = -26278.5  # N71: -26649.6006 -> -26278.5

$ws.Cells.Item(104, 8).Value = 0  # H104: 4000 -> 0
$ws.Cells.Item(104, 10).Value = 0  # J104: 4000 -> 0
$ws.Cells.Item(104, 12).Value = 0  # L104: 12000 -> 0
$ws.Cells.Item(104, 14).ClearContents()  # N104: -17242 -> (removed)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 48993.668  # H39: 48994.75 -> 48993.668
$ws.Cells.Item(39, 10).Value = 48993.668  # J39: 48994.75 -> 48993.668
$ws.Cells.Item(39, 12).Value = 48993.668  # L39: 48994.75 -> 48993.668
$ws.Cells.Item(39, 14).Value = -50057.668  # N39: -50058.75 -> -50057.668

$ws.Cells.Item(44, 8).Value = 5000  # H44: 0 -> 5000
$ws.Cells.Item(44, 10).Value = 5000  # J44: 0 -> 5000
$ws.Cells.Item(44, 12).Value = 5000  # L44: 0 -> 5000
$ws.Cells.Item(44, 14).Value = -6192  # N44: None -> -6192

$ws.Cells.Item(102, 8).Value = 3461.3076  # H102: 3763.6365 -> 3461.3076
$ws.Cells.Item(102, 9).Value = 2166.1667  # I102: 2350 -> 2166.1667
$ws.Cells.Item(102, 11).Value = 2166.1667  # K102: 2350 -> 2166.1667
$ws.Cells.Item(102, 13).Value = -544.1667000000002  # M102: -728 -> -544.1667000000002

$ws.Cells.Item(122, 8).Value = 3376.4375  # H122: 2933.8333 -> 3376.4375
$ws.Cells.Item(122, 9).Value = 2835.889  # I122: 2400.8125 -> 2835.889
$ws.Cells.Item(122, 10).Value = 4071.4285  # J122: 3999.875 -> 4071.4285
$ws.Cells.Item(122, 11).Value = 8507.667000000001  # K122: 7202.4375 -> 8507.667000000001
$ws.Cells.Item(122, 12).Value = 12214.2855  # L122: 11999.625 -> 12214.2855
$ws.Cells.Item(122, 13).Value = -6057.667000000001  # M122: -4752.4375 -> -6057.667000000001
$ws.Cells.Item(122, 14).Value = -17114.2855  # N122: -16899.625 -> -17114.2855

$ws.Cells.Item(126, 8).Value = 10495  # H126: 10613.941 -> 10495
$ws.Cells.Item(126, 9).Value = 19168.5  # I126: 17252.875 -> 19168.5
$ws.Cells.Item(126, 11).Value = 57505.5  # K126: 51758.625 -> 57505.5
$ws.Cells.Item(126, 13).Value = -55035.5  # M126: -49288.625 -> -55035.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 26319142  # H7: 27781122 -> 26319142
$ws.Cells.Item(7, 9).Value = 41669384  # I7: 38464250 -> 41669384
$ws.Cells.Item(7, 10).Value = 4443.7144  # J7: 4999.8 -> 4443.7144
$ws.Cells.Item(7, 11).Value = 41669384  # K7: 38464250 -> 41669384
$ws.Cells.Item(7, 12).Value = 4443.7144  # L7: 4999.8 -> 4443.7144
$ws.Cells.Item(7, 13).Value = -41669272  # M7: -38464138 -> -41669272
$ws.Cells.Item(7, 14).Value = -4667.7144  # N7: -5223.8 -> -4667.7144

$ws.Cells.Item(100, 8).Value = 11886.889  # H100: 9280.916999999999 -> 11886.889
$ws.Cells.Item(100, 10).Value = 0  # J100: 1463 -> 0
$ws.Cells.Item(100, 12).Value = 0  # L100: 1463 -> 0
$ws.Cells.Item(100, 14).ClearContents()  # N100: -2545 -> (removed)

$ws.Cells.Item(122, 8).Value = 4551.737  # H122: 4034.76 -> 4551.737
$ws.Cells.Item(122, 9).Value = 3900  # I122: 3299.0667 -> 3900
$ws.Cells.Item(122, 11).Value = 11700  # K122: 9897.2001 -> 11700
$ws.Cells.Item(122, 13).Value = -9250  # M122: -7447.2001 -> -9250

$ws.Cells.Item(126, 8).Value = 26319142  # H126: 27781122 -> 26319142
$ws.Cells.Item(126, 9).Value = 41669384  # I126: 38464250 -> 41669384
$ws.Cells.Item(126, 10).Value = 4443.7144  # J126: 4999.8 -> 4443.7144
$ws.Cells.Item(126, 11).Value = 125008152  # K126: 115392750 -> 125008152
$ws.Cells.Item(126, 12).Value = 13331.1432  # L126: 14999.4 -> 13331.1432
$ws.Cells.Item(126, 13).Value = -125005682  # M126: -115390280 -> -125005682
$ws.Cells.Item(126, 14).Value = -18271.1432  # N126: -19939.4 -> -18271.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2923.1667  # H122: 2979 -> 2923.1667
$ws.Cells.Item(122, 9).Value = 1985.125  # I122: 1996.25 -> 1985.125
$ws.Cells.Item(122, 10).Value = 4799.25  # J122: 5599.6665 -> 4799.25
$ws.Cells.Item(122, 11).Value = 5955.375  # K122: 5988.75 -> 5955.375
$ws.Cells.Item(122, 12).Value = 14397.75  # L122: 16798.9995 -> 14397.75
$ws.Cells.Item(122, 13).Value = -3505.375  # M122: -3538.75 -> -3505.375
$ws.Cells.Item(122, 14).Value = -19297.75  # N122: -21698.9995 -> -19297.75
